$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Auto-update data + news: refresh Revolving Consumer Credit row (row 6)
$ws.Range("E6").Value = 1313920.25

# Force the as_of cell to remain plain text (avoid Excel's "Nov 2025" -> date
# auto-conversion), then strip the temporary text format so the cell's style
# matches the rest of the column (no explicit style id).
$ws.Range("F6").NumberFormat = "@"
$ws.Range("F6").Value = "Nov 2025"
$ws.Range("F6").ClearFormats()

$ws.Range("G6").Value = 1097707.535371901
$ws.Range("H6").Value = -25154.75
$ws.Range("I6").Value = -0.01878516886656834
